# Append the new resale-numbers row (2025-01-03 22:40:41) as row 6,
# matching the existing table of Date/Time/Weekday/Week + 16 city columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A ("Date") and D ("Week") look numeric/date-like ("2025-01-03",
# "00") so Excel's smart-parsing would otherwise coerce them into a date
# serial / plain number. Force Text format before assigning, then drop the
# formatting back to Normal so the new row doesn't pick up a visible style
# (matches rows 2-5, which carry no explicit style either).
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2025-01-03"
$ws.Range("A6").Style = "Normal"

$ws.Range("B6").Value = "22:40:41"
$ws.Range("C6").Value = "Friday"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "00"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = 127115
$ws.Range("F6").Value = 143598
$ws.Range("G6").Value = 168124
$ws.Range("H6").Value = 158448
$ws.Range("I6").Value = -1
$ws.Range("J6").Value = 142022
$ws.Range("K6").Value = -1
$ws.Range("L6").Value = -1
$ws.Range("M6").Value = 192275
$ws.Range("N6").Value = 114626
$ws.Range("O6").Value = 45396
$ws.Range("P6").Value = 28276
$ws.Range("Q6").Value = 63218
$ws.Range("R6").Value = -1
$ws.Range("S6").Value = 48471
$ws.Range("T6").Value = -1
